# Update timestamps in the handback-status workbook (Generate Report for Handback)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (G2)
$wsOverview.Range("G2").Value = "2016-09-04 09:10:14"

# zh-cn sheet: Correspond Handoff Datetime (H2), Correspond Handback DateTime (K2)
$wsZhCn.Range("H2").Value = "2016-09-04 09:10:06"
$wsZhCn.Range("K2").Value = "2016-09-04 09:10:32"

# de-de sheet: Correspond Handoff Datetime (H2), Correspond Handback DateTime (K2)
$wsDeDe.Range("H2").Value = "2016-09-04 09:10:14"
$wsDeDe.Range("K2").Value = "2016-09-04 09:10:38"
